$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

# Pre-format the numeric-looking text cells as Text so Excel keeps them as
# strings instead of silently converting them to numbers on input.
$wsFollower.Range("B2").NumberFormat = "@"
$wsFollower.Range("D2").NumberFormat = "@"
$wsFollower.Range("E2").NumberFormat = "@"
$wsFollower.Range("F2").NumberFormat = "@"
$wsFollower.Range("B3").NumberFormat = "@"
$wsFollower.Range("D3").NumberFormat = "@"
$wsFollower.Range("E3").NumberFormat = "@"
$wsFollower.Range("F3").NumberFormat = "@"
$wsFollower.Range("B4").NumberFormat = "@"
$wsFollower.Range("D4").NumberFormat = "@"
$wsFollower.Range("E4").NumberFormat = "@"
$wsFollower.Range("F4").NumberFormat = "@"

$wsFollower.Range("A2").Value = "-26.483333333333338 - 3x + 6.88888888888889y"
$wsFollower.Range("B2").Value = "29.483333333333338"
$wsFollower.Range("D2").Value = "0.13"
$wsFollower.Range("E2").Value = "9.6"
$wsFollower.Range("F2").Value = "6.2"

$wsFollower.Range("A3").Value = "-1.6175000000000002 + x - 0.050000000000000044y"
$wsFollower.Range("B3").Value = "-2.3825"
$wsFollower.Range("D3").Value = "0.21"
$wsFollower.Range("E3").Value = "7.800000000000001"
$wsFollower.Range("F3").Value = "9.7"

$wsFollower.Range("A4").Value = "-9.049999999999999 + x + 0.6666666666666667y"
$wsFollower.Range("B4").Value = "-2.0499999999999994"
$wsFollower.Range("D4").Value = "0.91"
$wsFollower.Range("E4").Value = "8.7"
$wsFollower.Range("F4").Value = "0.6"

# --- Sheet: Punto_modificado ---
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2").NumberFormat = "@"
$wsPunto.Range("B2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "1.85"
$wsPunto.Range("B2").Value = "4.65"

# --- Sheet: Vector_bf ---
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and Worksheets.Item
# does case-insensitive name lookups (always matching the first sheet added,
# "Vector_bf") - so these two must be addressed by their 1-based tab index
# instead of by name to land on the correct sheet.
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2").NumberFormat = "@"
$wsBf.Range("A2").Value = "-6.0167222222222225"

# --- Sheet: Vector_BF ---
$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2").NumberFormat = "@"
$wsBF.Range("A3").NumberFormat = "@"
$wsBF.Range("A2").Value = "19.949999999999996"
$wsBF.Range("A3").Value = "-114.54333333333334"

# --- Sheet: Vector_Alpha ---
$wsAlpha = $wb.Worksheets.Item(7)
$wsAlpha.Range("A2").Value = 0.8999999999999999
